$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..41 (timestamp shifted by +1 day; Dragosel added to portfolio;
# the last two time-slots of the previous day are dropped).
$data = @(
    @(45890, 8.609999999999999, 0.157),
    @(45890.01041666666, 18.067, 0),
    @(45890.02083333334, 10.5, 0),
    @(45890.03125, 5.889, 0.155),
    @(45890.04166666666, 0.676, 1.479),
    @(45890.05208333334, 8.225, 0),
    @(45890.0625, 11.659, 0),
    @(45890.07291666666, 15.023, 0),
    @(45890.08333333334, 4.94, 0),
    @(45890.09375, 9.249000000000001, 0),
    @(45890.10416666666, 14.138, 0),
    @(45890.11458333334, 11.272, 0),
    @(45890.125, 19.01, 0),
    @(45890.13541666666, 6.605, 0),
    @(45890.14583333334, 8.164999999999999, 0),
    @(45890.15625, 5.749, 0),
    @(45890.16666666666, 0.8139999999999999, 0.489),
    @(45890.17708333334, 0.164, 3.46),
    @(45890.1875, 0, 10.436),
    @(45890.19791666666, 0.408, 2.51),
    @(45890.20833333334, 0.226, 5.963),
    @(45890.21875, 0.26, 2.21),
    @(45890.22916666666, 0.607, 0.802),
    @(45890.23958333334, 8.199999999999999, 0.008999999999999999),
    @(45890.25, 8.904999999999999, 0.038),
    @(45890.26041666666, 5.483, 0.008999999999999999),
    @(45890.27083333334, 2.083, 0.358),
    @(45890.28125, 10.778, 0.014),
    @(45890.29166666666, 3.548, 0.9340000000000001),
    @(45890.30208333334, 0.574, 0.8149999999999999),
    @(45890.3125, 9.824, 0.011),
    @(45890.32291666666, 0.12, 12.754),
    @(45890.33333333334, 6.006, 6.814),
    @(45890.34375, 11.335, 0),
    @(45890.35416666666, 14.689, 0),
    @(45890.36458333334, 1.621, 1.315),
    @(45890.375, 33.867, 0),
    @(45890.38541666666, 36.845, 0),
    @(45890.39583333334, 39.953, 0),
    @(45890.40625, 14.36, 0.021)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove the now-unused trailing rows (old rows 42 and 43)
$lastRow = $startRow + $data.Count
$ws.Range("A$($lastRow):C43").Clear()
